$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $value) {
    # Force text storage so date-like strings (e.g. "2011-11-25") are not
    # auto-converted to a date serial number, then restore the default
    # ("Normal") cell style so no stray number-format style lingers.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet "土地" (land): add the first data row.
# ---------------------------------------------------------------------
$wsLand = $wb.Worksheets.Item("土地")
$wsLand.Cells.Item(2, 1).Value = 14
$wsLand.Cells.Item(2, 2).Value = "高雄市美濃區福安段02410000地號"
$wsLand.Cells.Item(2, 3).Value = 1140.65
$wsLand.Cells.Item(2, 4).Value = "12分之1"
$wsLand.Cells.Item(2, 5).Value = "李永得"
$wsLand.Cells.Item(2, 6).Value = "95年04月10日"
$wsLand.Cells.Item(2, 7).Value = "貝賈"
$wsLand.Cells.Item(2, 8).Value = 133075.83
$wsLand.Cells.Item(2, 9).Value = "land"
$wsLand.Cells.Item(2, 10).Value = "normal"
Set-TextValue $wsLand.Cells.Item(2, 11) "2011-11-25"
$wsLand.Cells.Item(2, 12).Value = "邱議瑩"
$wsLand.Cells.Item(2, 13).Value = 913
$wsLand.Cells.Item(2, 14).Value = "tmp5431"
$wsLand.Cells.Item(2, 15).Value = 14
$wsLand.Cells.Item(2, 16).Value = 0.0833333333333333
$wsLand.Cells.Item(2, 17).Value = 95.0541666666667

# ---------------------------------------------------------------------
# Sheet "存款" (deposit): insert a new row 2 (clone of the header record)
# ahead of the existing data, shifting everything else down by one row.
# ---------------------------------------------------------------------
$wsDeposit = $wb.Worksheets.Item("存款")
$wsDeposit.Rows.Item(2).Insert()
$wsDeposit.Range("A2:G2").Style = $wsDeposit.Range("A3:G3").Style
$wsDeposit.Cells.Item(2, 1).Value = 18
$wsDeposit.Cells.Item(2, 2).Value = "臺灣銀行城中分行"
$wsDeposit.Cells.Item(2, 3).Value = "綜合存款"
$wsDeposit.Cells.Item(2, 4).Value = "美金"
$wsDeposit.Cells.Item(2, 5).Value = "邱議瑩"
$wsDeposit.Cells.Item(2, 6).Value = 2000
$wsDeposit.Cells.Item(2, 7).Value = 59440

# ---------------------------------------------------------------------
# Sheet "股票" (stock): insert a new row 2 for the Equinox Minerals stake,
# shifting the existing 中興商銀 row down to row 3.
# ---------------------------------------------------------------------
$wsStock = $wb.Worksheets.Item("股票")
$wsStock.Rows.Item(2).Insert()
$wsStock.Range("A2:N2").Style = $wsStock.Range("A3:N3").Style
$wsStock.Cells.Item(2, 1).Value = 43
$wsStock.Cells.Item(2, 2).Value = "EquinoxMineralsLimited"
$wsStock.Cells.Item(2, 3).Value = "李永得"
$wsStock.Cells.Item(2, 4).Value = 3497
$wsStock.Cells.Item(2, 5).Value = 5.95
$wsStock.Cells.Item(2, 6).Value = "澳幣"
$wsStock.Cells.Item(2, 7).Value = 624215
$wsStock.Cells.Item(2, 8).Value = "stock"
$wsStock.Cells.Item(2, 9).Value = "normal"
Set-TextValue $wsStock.Cells.Item(2, 10) "2011-11-25"
$wsStock.Cells.Item(2, 11).Value = "邱議瑩"
$wsStock.Cells.Item(2, 12).Value = 913
$wsStock.Cells.Item(2, 13).Value = "tmp5431"
$wsStock.Cells.Item(2, 14).Value = 43

# ---------------------------------------------------------------------
# Sheet "保險" (insurance): insert a new row 2 (clone of the header
# record), shifting the existing 台灣人壽新祥和定期壽險 row down to row 3.
# ---------------------------------------------------------------------
$wsInsurance = $wb.Worksheets.Item("保險")
$wsInsurance.Rows.Item(2).Insert()
$wsInsurance.Range("A2:D2").Style = $wsInsurance.Range("A3:D3").Style
$wsInsurance.Cells.Item(2, 1).Value = 49
$wsInsurance.Cells.Item(2, 2).Value = "台灣人壽"
$wsInsurance.Cells.Item(2, 3).Value = "台灣人壽歲歲長泰還本终身險"
$wsInsurance.Cells.Item(2, 4).Value = "李永得"

# ---------------------------------------------------------------------
# Sheet "事業投資" (business investment): insert a new row 2 (clone of the
# header record for 玉山社事業股份有限公司), shifting the existing two
# rows down by one.
# ---------------------------------------------------------------------
$wsInvestment = $wb.Worksheets.Item("事業投資")
$wsInvestment.Rows.Item(2).Insert()
$wsInvestment.Range("A2:G2").Style = $wsInvestment.Range("A3:G3").Style
$wsInvestment.Cells.Item(2, 1).Value = 55
$wsInvestment.Cells.Item(2, 2).Value = "李永得"
$wsInvestment.Cells.Item(2, 3).Value = "玉山社事業股份有限公司"
$wsInvestment.Cells.Item(2, 4).Value = "臺北市大安區仁愛路四段145號3樓之2"
$wsInvestment.Cells.Item(2, 5).Value = 1000000
$wsInvestment.Cells.Item(2, 6).Value = "84年07月08日"
$wsInvestment.Cells.Item(2, 7).Value = "發起設立"
